# Weekly cryptos-list refresh (GitHub Actions scrape).
#
# For every changed row this updates the Price (column D) and/or the
# Volume(1h) (column E) cell to the newly scraped reading. Rows 12/13 and
# 48/49 additionally swap which coin (Coin/Link in columns B/C) occupies
# that rank, since WrappedEther/Polkadot and Mantle/EnergySwap traded
# places in the ranking that run.
#
# Column D holds plain-text numbers (e.g. "215.21", "25.785.61" using the
# sheet's "." thousands separator) rather than real numerics, so after
# writing each Price value we force the cell back to General/no-style via
# NumberFormat "@" + ClearFormats(); otherwise the COM layer would helpfully
# "recognize" strings like "215.21" as a genuine number and coerce the type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    # Write $value into $cell while guaranteeing it is stored as text,
    # matching the original inline-string cells (no numeric coercion,
    # no residual cell style).
    $range = $ws.Range($cell)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "25.785.61"
$ws.Range("E2").Value = "  -0.15%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.635.42"
$ws.Range("E3").Value = "  -0.06%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.00%  "

# Row 5 - BNB
Set-TextValue "D5" "215.21"
$ws.Range("E5").Value = "  -0.30%  "

# Row 6 - XRP
$ws.Range("E6").Value = "  -0.57%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.02%  "

# Row 8 - Cardano
Set-TextValue "D8" "0.257"
$ws.Range("E8").Value = "  -0.21%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.0641"
$ws.Range("E9").Value = "  -0.38%  "

# Row 10 - Solana
Set-TextValue "D10" "19.82"
$ws.Range("E10").Value = "  +1.10%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +1.02%  "

# Row 12 - WrappedEther
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D12" "1.645.93"
$ws.Range("E12").Value = "  +0.58%  "

# Row 13 - Polkadot
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D13" "4.24"
$ws.Range("E13").Value = "  -0.82%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue "D14" "1.862.26"
$ws.Range("E14").Value = "  -0.06%  "

# Row 15 - Polygon
Set-TextValue "D15" "0.556"
$ws.Range("E15").Value = "  -1.45%  "

# Row 16 - ShibaInu
Set-TextValue "D16" "0.0₃0775"
$ws.Range("E16").Value = "  +1.91%  "

# Row 17 - Litecoin
Set-TextValue "D17" "63.02"
$ws.Range("E17").Value = "  -0.31%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "25.804.57"
$ws.Range("E18").Value = "  -0.18%  "

# Row 20 - Uniswap
Set-TextValue "D20" "4.44"
$ws.Range("E20").Value = "  +2.65%  "

# Row 21 - BitcoinCash
Set-TextValue "D21" "193.75"
$ws.Range("E21").Value = "  -0.96%  "

# Row 22 - Avalanche
Set-TextValue "D22" "9.93"
$ws.Range("E22").Value = "  +0.36%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  +0.98%  "

# Row 24 - BinanceUSD
$ws.Range("E24").Value = "  -0.01%  "

# Row 25 - Toncoin
Set-TextValue "D25" "1.76"
$ws.Range("E25").Value = "  -1.66%  "

# Row 26 - Monero
Set-TextValue "D26" "139.63"
$ws.Range("E26").Value = "  -0.10%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -4.81%  "

# Row 28 - Cosmos
Set-TextValue "D28" "6.82"
$ws.Range("E28").Value = "  +0.48%  "

# Row 29 - EthereumClassic
Set-TextValue "D29" "15.54"
$ws.Range("E29").Value = "  +0.19%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.11%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  +1.34%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  +1.08%  "

# Row 33 - Filecoin
Set-TextValue "D33" "3.25"
$ws.Range("E33").Value = "  +0.60%  "

# Row 34 - LidoDAOToken
Set-TextValue "D34" "1.59"
$ws.Range("E34").Value = "  +2.27%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  +0.69%  "

# Row 36 - ARBITRUM
Set-TextValue "D36" "0.896"
$ws.Range("E36").Value = "  -0.90%  "

# Row 37 - MXToken
$ws.Range("E37").Value = "  -0.21%  "

# Row 38 - ImmutableX
Set-TextValue "D38" "0.549"
$ws.Range("E38").Value = "  -0.30%  "

# Row 39 - Maker
Set-TextValue "D39" "1.107.03"
$ws.Range("E39").Value = "  -2.02%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  +0.24%  "

# Row 41 - PaxDollar
$ws.Range("E41").Value = "  +0.62%  "

# Row 43 - Quant
Set-TextValue "D43" "99.20"
$ws.Range("E43").Value = "  +1.48%  "

# Row 44 - TrustWalletToken
Set-TextValue "D44" "0.800"
$ws.Range("E44").Value = "  -0.02%  "

# Row 45 - BabyDogeCoin
Set-TextValue "D45" "0.0₆0108"
$ws.Range("E45").Value = "  -4.17%  "

# Row 46 - SynthetixNetwork
Set-TextValue "D46" "2.56"
$ws.Range("E46").Value = "  +14.36%  "

# Row 47 - Aave
Set-TextValue "D47" "55.62"
$ws.Range("E47").Value = "  +0.40%  "

# Row 48 - Mantle
$ws.Range("B48").Value = "Mantle"
$ws.Range("C48").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D48" "0.418"
$ws.Range("E48").Value = "  -5.96%  "

# Row 49 - EnergySwap
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D49" "7.67"
$ws.Range("E49").Value = "  -0.08%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  -0.34%  "

# Row 51 - Frax
$ws.Range("E51").Value = "  -0.01%  "
